$wb = $excel.ActiveWorkbook

# --- Answers sheet: replace student DDD2/ffe with XXX111/Sami Alfattani ---
$ans = $wb.Worksheets.Item("Answers")
$ans.Range("A2").Value = "XXX111"
$ans.Range("B2").Value = "Sami Alfattani"
$ans.Range("A2:B2").Select()

# --- Timer sheet: same student swap + new recorded timing values ---
$timer = $wb.Worksheets.Item("Timer")
$timer.Range("A2").Value = "XXX111"
$timer.Range("B2").Value = "Sami Alfattani"
$timer.Range("C2:H2").NumberFormat = "@"
$timer.Range("C2").Value = "01:00"
$timer.Range("D2").Value = "01:00"
$timer.Range("E2").Value = "01:00"
$timer.Range("F2").Value = "01:00"
$timer.Range("G2").Value = "01:30"
$timer.Range("H2").Value = "05:30"
$timer.Range("H3").Select()

# Timer becomes the active sheet/tab (was BlackWhite_List before)
$timer.Activate()
